# Jeannine's log - add a new log entry (spell checked and permission changed)
#
# Appends a new "section separator" row + a new data row to the bottom of
# the "Logs" sheet, mirroring the existing pattern used throughout the
# sheet (a shaded separator row naming the day-of-week, immediately
# followed by the actual log entry row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# The most recent existing block lives in rows 29:30 (separator + entry).
# Duplicate it (values + formatting) into rows 34:35 for the new entry, then
# update only the date (and corresponding day-of-week label).
$src = $ws.Range("A29:F30")
$dst = $ws.Range("A34:F35")
$src.Copy($dst)

# New entry: Friday Aug 12, 2016 (serial 42594) - same task as before.
$ws.Cells.Item(34, 2).Value = "FRIDAY"
$ws.Cells.Item(35, 2).Value = 42594

# Preserve the wrapped-text row height used by every other detail row.
$ws.Rows.Item(35).RowHeight = 45

# Move the active selection to reflect where the user would continue typing.
$ws.Range("F39").Select()
